# Update ResultsTable metadata descriptions for HostUse_Combined (row 26) and
# Voltinism_Combined (row 27), and move the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 (A27 = "Voltinism_Combined"): reword the combined-categories description.
$ws.Range("B27").Value = 'Voltinism combined into a categorical variable two levels: 1) Semi-and univoltine species and 2) Multivoltine species (including facultative bivoltine species)'

# Row 26 (A26 = "HostUse_Combined"): reword the combined-categories description.
$ws.Range("B26").Value = 'HostUse combined into a categorical variable with two levels: 1) Specialists= including monophagous species feeding only on a single host plant species,  oligophagous species feeding on a limited number (two to four) of host plant species, and species that feed on lichen or fungi, and 2) Generalist= polyphagous species feeding on several (> 5) host plants.'

# Move the saved selection / active cell to H35, matching the author's last click.
$ws.Range("H35").Select()
